$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ===========================================================================
# Sheet1 ("Danh sach")
# ===========================================================================

# Updated file-creation date note (3/5/2023 -> 4/5/2023)
$ws1.Range("D4").Value = "Ngày tạo file: ngày 4 tháng 5, 2023"

# Section banner text change, and a new "Bai thi" label box next to it
$ws1.Range("B6").Value = "Thông tin thí sinh"
$ws1.Range("G6").Value = "Bài thi"
$ws1.Range("G6:I6").Merge()

# New header label for the subject/score column (note the trailing space)
$ws1.Range("I7").Value = "Test "

# Move the per-student score out of column G and into the new column I
# (under the "Test " header).
$ws1.Range("I8").Value = 0
$ws1.Range("I9").Value = "Chưa làm"
$ws1.Range("I10").Value = 0
$ws1.Range("I11").Value = "Chưa làm"

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------

# New bordered / coloured "Bai thi" box spanning G6:I6
$box = $ws1.Range("G6:I6")
$box.Borders.Color = 9779238
$box.Borders.LineStyle = 1
$box.Font.Color = 9779238
$box.Font.Bold = $true
$box.Font.Size = 10
$box.HorizontalAlignment = -4108

# Header row (row 7) - left align, vertically centered
$hdr = $ws1.Range("A7:I7")
$hdr.HorizontalAlignment = -4131
$hdr.VerticalAlignment = -4108

# Data rows (8-11) - left align
$data = $ws1.Range("A8:I11")
$data.HorizontalAlignment = -4131

# Now that the formatting sweep above has touched the full A8:I11 rectangle,
# drop the cells that should no longer exist now that the score moved to
# column I.
$ws1.Range("G8").Clear()
$ws1.Range("H8").Clear()
$ws1.Range("G9").Clear()
$ws1.Range("H9").Clear()
$ws1.Range("G11").Clear()
$ws1.Range("H11").Clear()

# The numeric "0" score cells get a bold red sz-9 font (highlights the
# ungraded / zero scores).
$ws1.Range("I8").Font.Bold = $true
$ws1.Range("I8").Font.Color = 2368716
$ws1.Range("I8").Font.Size = 9

$ws1.Range("I10").Font.Bold = $true
$ws1.Range("I10").Font.Color = 2368716
$ws1.Range("I10").Font.Size = 9

Write-Output "sheet1 done"

# ===========================================================================
# Sheet2 ("So con diem") - grading-scale histogram table
# ===========================================================================

$labels = @("<=0", "<=1", "<=2", "<=3", "<=4", "<=5", "<=6", "<=7", "<=8", "<=9", "<=10")
$counts = @(2, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 2).Value = $labels[$i]
    $ws2.Cells.Item($row, 3).Value = $counts[$i]
}

Write-Output "sheet2 done"
